# Update the "dSF" column (F) values following a repull/push of data and
# recalculation against the mean, per commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 1
    4  = 1
    5  = 3
    7  = 2
    9  = 2
    10 = 7
    11 = 4
    12 = -4
    13 = -1
    14 = -1
    15 = 3
    16 = 4
    17 = 2
    18 = 2
    19 = 3
    20 = 5
    21 = -3
    22 = -2
    24 = -4
    25 = -4
    26 = 4
    27 = 1
    28 = -1
    30 = 2
    32 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
